{"js": "// Bump the sprint number (1 -> 2) and update the review date\n// (02/09/18 -> 02/21/18) in the document-info table at the top of the\n// checklist. Cells are located by their known row position relative to\n// the \"Sprint No.\" / \"Review Date\" label cells so the script still finds\n// the right targets even if unrelated rows are inserted/removed upstream.\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst infoTable = tables.items[0];\nconst rows = infoTable.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    cell.load(\"value\");\n  }\n}\nawait context.sync();\n\nlet sprintCell = null;\nlet dateCell = null;\n\nfor (const row of rows.items) {\n  const cells = row.cells.items;\n  for (let i = 0; i < cells.length - 1; i++) {\n    const label = (cells[i].value || \"\").trim();\n    if (label === \"Sprint No.\") {\n      sprintCell = cells[i + 1];\n    } else if (label === \"Review Date\") {\n      dateCell = cells[i + 1];\n    }\n  }\n}\n\nif (!sprintCell || !dateCell) {\n  throw new Error(\"Could not locate Sprint No. / Review Date cells\");\n}\n\nsprintCell.getRange().insertText(\"2\", \"Replace\");\ndateCell.getRange().insertText(\"02/21/18\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Bump the sprint number (1 -> 2) and update the review date\n# (02/09/18 -> 02/21/18) in the document-info table at the top of the\n# checklist. Cells are located by their known label text (\"Sprint No.\" /\n# \"Review Date\") so the script still finds the right targets even if\n# unrelated rows are inserted/removed upstream.\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$sprintCell = $null\n$dateCell = $null\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $row = $t.Rows($r)\n    $n = $row.Cells.Count\n    for ($c = 1; $c -lt $n; $c++) {\n        $label = $row.Cells($c).Range.Text.TrimEnd([char]13, [char]7).Trim()\n        if ($label -eq \"Sprint No.\") {\n            $sprintCell = $row.Cells($c + 1)\n        } elseif ($label -eq \"Review Date\") {\n            $dateCell = $row.Cells($c + 1)\n        }\n    }\n}\n\n$sprintCell.Range.Text = \"2\"\n$dateCell.Range.Text = \"02/21/18\"\n"}
